# fix subset method definition, remove unneeded plotSpec() override
#
# Adds a new "components" worksheet (mirroring the existing fGroups /
# mslists / formulas / compounds sheets) right after "compounds", fills it
# with the set-method implementation-status table, and makes it the active
# (selected) sheet - matching the selection change left behind on the
# "compounds" sheet as well.

$wb = $excel.ActiveWorkbook

# The "compounds" sheet was the active/selected one before this edit
# (tabSelected) with G28 selected; afterwards it keeps its data but the
# selection moves to the header row and the tab selection moves to the new
# sheet. Set that selection now, while "compounds" is still the active
# sheet, so it persists correctly once the new sheet is activated.
$compounds = $wb.Worksheets.Item("compounds")
$compounds.Range("B1:G1").Select()

# Insert the new sheet right after "compounds" and rename it.
$ws = $wb.Worksheets.Add([System.Type]::Missing, $compounds)
$ws.Name = "components"

# Column width tweak (matches the narrower "A" column used on this sheet).
$ws.Columns.Item(1).ColumnWidth = 16.140625

# Header row.
$ws.Range("B1").Value = "as-is"
$ws.Range("C1").Value = "almost as-is"
$ws.Range("D1").Value = "implement"
$ws.Range("E1").Value = "not supported"
$ws.Range("F1").Value = "ionize"
$ws.Range("G1").Value = "done"

# Data rows - one row per `components` S4 method/accessor, with X marks in
# the column describing its current set-support status.
$ws.Range("A2").Value = '$'
$ws.Range("B2").Value = "X"
$ws.Range("G2").Value = "X"

$ws.Range("A3").Value = "["
$ws.Range("C3").Value = "X"
$ws.Range("G3").Value = "X"

$ws.Range("A4").Value = "[["
$ws.Range("B4").Value = "X"
$ws.Range("G4").Value = "X"

$ws.Range("A5").Value = "as.data.table"
$ws.Range("B5").Value = "X"
$ws.Range("G5").Value = "X"

$ws.Range("A6").Value = "componentInfo"
$ws.Range("B6").Value = "X"
$ws.Range("G6").Value = "X"

$ws.Range("A7").Value = "componentTable"
$ws.Range("B7").Value = "X"
$ws.Range("G7").Value = "X"

$ws.Range("A8").Value = "consensus"
$ws.Range("E8").Value = "X"

$ws.Range("A9").Value = "filter"
$ws.Range("C9").Value = "X"
$ws.Range("G9").Value = "X"

$ws.Range("A10").Value = "findFGroup"
$ws.Range("B10").Value = "X"
$ws.Range("G10").Value = "X"

$ws.Range("A11").Value = "groupNames"
$ws.Range("B11").Value = "X"
$ws.Range("G11").Value = "X"

$ws.Range("A12").Value = "initialize"
$ws.Range("C12").Value = "X"
$ws.Range("G12").Value = "X"

$ws.Range("A13").Value = "length"
$ws.Range("B13").Value = "X"
$ws.Range("G13").Value = "X"

$ws.Range("A14").Value = "names"
$ws.Range("B14").Value = "X"
$ws.Range("G14").Value = "X"

$ws.Range("A15").Value = "plotEIC"
$ws.Range("B15").Value = "X"
$ws.Range("D15").Value = "X"
$ws.Range("G15").Value = "X"
$ws.Range("H15").Value = "Seems enough, assuming we're not planning to merge components"

$ws.Range("A16").Value = "plotEICHash"
$ws.Range("B16").Value = "X"
$ws.Range("G16").Value = "X"

$ws.Range("A17").Value = "plotSpec"
$ws.Range("B17").Value = "X"
$ws.Range("D17").Value = "X"
$ws.Range("G17").Value = "X"
$ws.Range("H17").Value = "Seems enough, assuming we're not planning to merge components"

$ws.Range("A18").Value = "plotSpecHash"
$ws.Range("B18").Value = "X"
$ws.Range("G18").Value = "X"

$ws.Range("A19").Value = "show"
$ws.Range("C19").Value = "X"
$ws.Range("G19").Value = "X"

# Match the author's final cursor position / selection on the new sheet.
$ws.Range("H17").Select()
